{"js": "// The document's credits line reads:\n//   \"created by Ola Wiebe, Malena Bengtsson and Victor Lastname 2019\"\n// The placeholder surname \"Lastname\" needs to become \"T\u00f6rnbom\" (the\n// co-author's real last name), as reflected by the canonical OOXML diff\n// (the run that used to contain \"Lastname\" now contains \"T\u00f6rnbom\").\nconst body = context.document.body;\n\nconst results = body.search(\"Lastname\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"T\u00f6rnbom\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document's credits line reads:\n#   \"created by Ola Wiebe, Malena Bengtsson and Victor Lastname 2019\"\n# The placeholder surname \"Lastname\" needs to become \"T\u00f6rnbom\" (the\n# co-author's real last name), as reflected by the canonical OOXML diff\n# (the run that used to contain \"Lastname\" now contains \"T\u00f6rnbom\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Lastname\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"T\u00f6rnbom\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\"Lastname\", $true, $true, $false, $false, $false, $true, 1, $false, \"T\u00f6rnbom\", 2)\n"}
